# Update the "Lista de Vinos" worksheet with the new wine ranking data
# and widen the B:D columns to match the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lista de Vinos")

# New data rows (score, wine, winery, varietals, price), replacing the
# previous ranking/content of rows 2-11.
$data = @(
    @(97,    "Finca Altamira",         "Bodega Norton",    '["Malbec Reserva - Tipo Uva: Malbec","Cabernet Sauvignon Gran Reserva - Tipo Uva: Cabernet Sauvignon"]', 7200),
    @(94.5,  "Estrellas Blend",        "Bodega El Esteco", '["Malbec Reserva - Tipo Uva: Malbec","Syrah Premium - Tipo Uva: Syrah"]', 6500),
    @(80.5,  "Cumbres Malbec",         "Bodega Norton",    '["Malbec Reserva - Tipo Uva: Malbec"]', 4800),
    @(63,    "Viña de las Estrellas",  "Bodega Norton",    '["Malbec Reserva - Tipo Uva: Malbec"]', 8000),
    @(61.43, "Trapiche Merlot",        "Bodega Trapiche",  '["Merlot Reserva Especial - Tipo Uva: Merlot"]', 3200),
    @(60.67, "Norton Malbec Especial", "Bodega Norton",    '["Malbec Reserva - Tipo Uva: Malbec"]', 5500),
    @(60,    "Cabernet de Trapiche",   "Bodega Trapiche",  '["Cabernet Sauvignon Gran Reserva - Tipo Uva: Cabernet Sauvignon"]', 4300),
    @(58.6,  "Blend de Valle",         "Bodega Norton",    '["Malbec Reserva - Tipo Uva: Malbec","Merlot Reserva Especial - Tipo Uva: Merlot"]', 6500),
    @(50.67, "Norton Reserva",         "Bodega Norton",    '["Malbec Reserva - Tipo Uva: Malbec"]', 7500),
    @(49.33, "Altura Malbec",          "Bodega Norton",    '["Malbec Reserva - Tipo Uva: Malbec"]', 4200)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
}

# Widen columns B, C, D; keep A and E at their original width.
# (9.1667 / 19.1667 are the ColumnWidth inputs that round-trip through
# Excel's internal pixel-based storage to the exact stored widths of
# 10 and 20 respectively.)
$ws.Range("A1").ColumnWidth = 9.1667
$ws.Range("B1:D1").ColumnWidth = 19.1667
$ws.Range("E1").ColumnWidth = 9.1667
